$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh Price (D) and Volume(1h) (E) columns with the latest crypto snapshot.
# Some Price values are plain decimals that Excel would otherwise auto-convert
# to numbers; format those cells as Text first so they stay literal strings,
# matching the source feed formatting used throughout column D.

$ws.Range("D2").Value = "66.168.05"
$ws.Range("E2").Value = "  +7.28%  "
$ws.Range("D3").Value = "3.015.63"
$ws.Range("E3").Value = "  +4.41%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.00"
$ws.Range("E5").Value = "  +3.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.94"
$ws.Range("E6").Value = "  +9.11%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "3.011.54"
$ws.Range("E8").Value = "  +4.32%  "
$ws.Range("E9").Value = "  +3.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.95"
$ws.Range("E10").Value = "  +1.08%  "
$ws.Range("E11").Value = "  +6.71%  "
$ws.Range("E12").Value = "  +5.81%  "
$ws.Range("E13").Value = "  +8.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.57"
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("D16").Value = "66.149.37"
$ws.Range("E16").Value = "  +7.26%  "
$ws.Range("D17").Value = "3.514.98"
$ws.Range("E17").Value = "  +4.38%  "
$ws.Range("E18").Value = "  +7.05%  "
$ws.Range("D19").Value = "3.015.00"
$ws.Range("E19").Value = "  +4.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "462.84"
$ws.Range("E20").Value = "  +7.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.89"
$ws.Range("E21").Value = "  +6.95%  "
$ws.Range("E22").Value = "  +5.54%  "
$ws.Range("E23").Value = "  +8.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.18"
$ws.Range("E24").Value = "  +4.19%  "
$ws.Range("E25").Value = "  +13.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.48"
$ws.Range("E26").Value = "  +5.04%  "
$ws.Range("E27").Value = "  +7.95%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.05"
$ws.Range("E29").Value = "  +14.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.36"
$ws.Range("E30").Value = "  +16.76%  "
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.61"
$ws.Range("E32").Value = "  +5.24%  "
$ws.Range("E33").Value = "  +5.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.02"
$ws.Range("E34").Value = "  +6.36%  "
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.996"
$ws.Range("E36").Value = "  +4.20%  "
$ws.Range("E37").Value = "  +8.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.16"
$ws.Range("E38").Value = "  +13.19%  "
$ws.Range("E39").Value = "  +9.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.49"
$ws.Range("E40").Value = "  +1.45%  "
$ws.Range("E41").Value = "  +8.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.304"
$ws.Range("E42").Value = "  +14.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "43.49"
$ws.Range("E43").Value = "  +10.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.43"
$ws.Range("E44").Value = "  +3.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "393.09"
$ws.Range("E45").Value = "  +14.50%  "
$ws.Range("D46").Value = "2.803.51"
$ws.Range("E46").Value = "  +4.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0354"
$ws.Range("E47").Value = "  +5.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.73"
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.59"
$ws.Range("E50").Value = "  +9.71%  "
$ws.Range("E51").Value = "  +4.42%  "
